# Updated cryptos list on Tue Oct  1 11:27:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price-column (D) values must stay plain text (the source sheet stores
# them as inline strings, not numbers), so force text format before writing
# any of them to avoid Excel auto-converting numeric-looking strings and
# dropping formatting such as trailing zeros.

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.872.56"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.94"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.85"
$ws.Range("E5").Value = "  +0.84%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.90"
$ws.Range("E6").Value = "  +1.30%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.34%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.11%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.43%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.79%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +0.21%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.94%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.77"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.110.58"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +1.07%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.799.36"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.671.58"
$ws.Range("E17").Value = "  +2.42%  "

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.20"
$ws.Range("E18").Value = "  +1.23%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("E19").Value = "  +4.00%  "

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").Value = "  -2.05%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.32"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.01%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.39"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24 - SuiNetwork
$ws.Range("E24").Value = "  +9.24%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +4.29%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +4.60%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28 - Bittensor
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "582.35"
$ws.Range("E28").Value = "  +0.43%  "

# Row 29 - Aptos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.32"
$ws.Range("E29").Value = "  +5.32%  "

# Row 30 - was Kaspa, now Binance-PegBSC-USD
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.03"
$ws.Range("E30").Value = "  +2.86%  "

# Row 31 - was Binance-PegBSC-USD, now Kaspa
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.161"
$ws.Range("E31").Value = "  -0.01%  "

# Row 32 - PancakeSwap
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  +0.79%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  +2.77%  "

# Row 34 - RenderToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.68"
$ws.Range("E34").Value = "  +3.21%  "

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +3.35%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  -0.95%  "

# Row 37 - EthereumClassic
$ws.Range("E37").Value = "  -0.15%  "

# Row 38 - FirstDigitalUSD
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.12%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +2.96%  "

# Row 40 - Monero
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.22"
$ws.Range("E40").Value = "  +0.56%  "

# Row 41 - dogwifhat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.57"
$ws.Range("E41").Value = "  +8.71%  "

# Row 42 - USDe
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "163.42"
$ws.Range("E43").Value = "  +4.65%  "

# Row 44 - InjectiveProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.28"
$ws.Range("E44").Value = "  +6.09%  "

# Row 45 - Filecoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.93"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46 - Hedera
$ws.Range("E46").Value = "  -0.52%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +0.91%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -0.75%  "

# Row 49 - VeChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0249"
$ws.Range("E49").Value = "  -0.53%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0238"
$ws.Range("E50").Value = "  +2.95%  "

# Row 51 - ONDO
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.790"
$ws.Range("E51").Value = "  +2.71%  "
